$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "AMSIN" (1st sheet): append two new regression-history rows
# (121, 180lastrun / 122, 180fnlrun) and correct the precise run-time
# timestamp that had been truncated on row 120.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("AMSIN")

$ws1.Cells.Item(120, 2).Value = 45089.67016552083

$ws1.Range("A119:G120").Copy()
$ws1.Range("A120:G122").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Cells.Item(121, 1).Value = "2023-07-31"
$ws1.Cells.Item(121, 2).Value = 45138.37467510417
$ws1.Cells.Item(121, 3).Value = "180lastrun"
$ws1.Cells.Item(121, 4).Value = 269
$ws1.Cells.Item(121, 5).Value = 262
$ws1.Cells.Item(121, 6).Value = 7
$ws1.Cells.Item(121, 7).Value = 6.52

$ws1.Cells.Item(122, 1).Value = "2023-07-31"
$ws1.Cells.Item(122, 2).Value = 45138.3927577662
$ws1.Cells.Item(122, 3).Value = "180fnlrun"
$ws1.Cells.Item(122, 4).Value = 269
$ws1.Cells.Item(122, 5).Value = 269
$ws1.Cells.Item(122, 6).Value = 0
$ws1.Cells.Item(122, 7).Value = 4.04

# ---------------------------------------------------------------------
# Sheet "BETA" (2nd sheet): append the 180beta run
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("BETA")

$ws2.Range("A38:G38").Copy()
$ws2.Range("A39:G39").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Cells.Item(39, 1).Value = "2023-08-01"
$ws2.Cells.Item(39, 2).Value = 45139.52310601852
$ws2.Cells.Item(39, 3).Value = "180beta"
$ws2.Cells.Item(39, 4).Value = 269
$ws2.Cells.Item(39, 5).Value = 259
$ws2.Cells.Item(39, 6).Value = 10
$ws2.Cells.Item(39, 7).Value = 8.22

# ---------------------------------------------------------------------
# Sheet "AMS" (3rd sheet): append the two HTFX runs (styled, copied
# forward from the previous row) plus the freshly typed "180live" row.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("AMS")

$ws3.Range("A81:G81").Copy()
$ws3.Range("A82:G83").PasteSpecial(-4122)
$ws3.Application.CutCopyMode = $false

$ws3.Cells.Item(82, 1).Value = "2023-07-18"
$ws3.Cells.Item(82, 2).Value = 45125.7645241088
$ws3.Cells.Item(82, 3).Value = "179htfxslots"
$ws3.Cells.Item(82, 4).Value = 269
$ws3.Cells.Item(82, 5).Value = 269
$ws3.Cells.Item(82, 6).Value = 0
$ws3.Cells.Item(82, 7).Value = 3.77

$ws3.Cells.Item(83, 1).Value = "2023-07-31"
$ws3.Cells.Item(83, 2).Value = 45138.81799356481
$ws3.Cells.Item(83, 3).Value = "179scndhtfx"
$ws3.Cells.Item(83, 4).Value = 269
$ws3.Cells.Item(83, 5).Value = 266
$ws3.Cells.Item(83, 6).Value = 3
$ws3.Cells.Item(83, 7).Value = 4.68

$ws3.Cells.Item(84, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(84, 1).Value = "2023-08-01"
$ws3.Cells.Item(84, 2).Value = 45139.83888939953
$ws3.Cells.Item(84, 3).Value = "180live"
$ws3.Cells.Item(84, 4).Value = 269
$ws3.Cells.Item(84, 5).Value = 265
$ws3.Cells.Item(84, 6).Value = 4
$ws3.Cells.Item(84, 7).Value = 5.5
